$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.4391767347414131
$ws0.Range("C2").Value = 0.0291736642814527
$ws0.Range("B3").Value = 0.1500138330350813
$ws0.Range("C3").Value = 0.8120492616816498
$ws0.Range("B4").Value = 0.1971588882678229
$ws0.Range("C4").Value = 0.7493492955355091

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -1.231142294725017
$ws1.Range("C2").Value = -0.6286683659443316
$ws1.Range("B3").Value = 0.158834666381007
$ws1.Range("C3").Value = 0.8549697779871629
$ws1.Range("B4").Value = 0.3343977303686733
$ws1.Range("C4").Value = 0.3590014067971175
